# Updated legacy GSC export data:
# The oldest day of data (2025-10-10, the first data row) has aged out of the
# export window, so remove that row from the "Chart" sheet. Excel shifts the
# remaining rows up, drops the now-unused shared string, and shrinks the
# used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the 2025-10-10 data (row 1 is the header row).
$ws.Rows.Item(2).Delete()
